$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "LaMelo Ball"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Charlotte Hornets"
$ws.Range("A3").Value = "Jusuf Nurkic"
$ws.Range("B3").Value = "C"
$ws.Range("C3").Value = "Charlotte Hornets"
$ws.Range("A4").Value = "Coby White"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("C4").Value = "Chicago Bulls"
$ws.Range("A5").Value = "Santi Aldama"
$ws.Range("B5").Value = "PF,C"
$ws.Range("C5").Value = "Memphis Grizzlies"
$ws.Range("A6").Value = "Precious Achiuwa"
$ws.Range("B6").Value = "PF,C"
$ws.Range("C6").Value = "New York Knicks"
$ws.Range("A7").Value = "Isaiah Hartenstein"
$ws.Range("B7").Value = "C"
$ws.Range("C7").Value = "Oklahoma City Thunder"
$ws.Range("A8").Value = "Naz Reid"
$ws.Range("B8").Value = "PF,C"
$ws.Range("C8").Value = "Minnesota Timberwolves"
$ws.Range("A9").Value = "Ausar Thompson"
$ws.Range("B9").Value = "SF,PF"
$ws.Range("C9").Value = "Detroit Pistons"
$ws.Range("A10").Value = "Malik Monk"
$ws.Range("B10").Value = "PG,SG,SF"
$ws.Range("C10").Value = "Sacramento Kings"
$ws.Range("A11").Value = "Derrick White"
$ws.Range("B11").Value = "PG,SG"
$ws.Range("C11").Value = "Boston Celtics"
$ws.Range("A12").Value = "Carlton Carrington"
$ws.Range("B12").Value = "PG,SG"
$ws.Range("C12").Value = "Washington Wizards"
$ws.Range("A13").Value = "Cade Cunningham"
$ws.Range("B13").Value = "PG,SG"
$ws.Range("C13").Value = "Detroit Pistons"
$ws.Range("A14").Value = "Onyeka Okongwu"
$ws.Range("B14").Value = "PF,C"
$ws.Range("C14").Value = "Atlanta Hawks"
$ws.Range("A15").Value = "Devin Vassell"
$ws.Range("B15").Value = "SG,SF"
$ws.Range("C15").Value = "San Antonio Spurs"
$ws.Range("A16").Value = "Damian Lillard"
$ws.Range("B16").Value = "PG"
$ws.Range("C16").Value = "Milwaukee Bucks"
$ws.Range("A17").Value = "Collin Sexton"
$ws.Range("B17").Value = "PG,SG"
$ws.Range("C17").Value = "Utah Jazz"
$ws.Range("A18").Value = "Anthony Davis"
$ws.Range("B18").Value = "PF,C"
$ws.Range("C18").Value = "Dallas Mavericks"
$ws.Range("A19").Value = "Andrew Wiggins"
$ws.Range("B19").Value = "SF,PF"
$ws.Range("C19").Value = "Miami Heat"
